# LacI model parameter boundaries: narrow down parameter ranges after
# manual parameter estimation. Also update the window view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the B (bmin) and C (bmax) columns -----------------------------
# Column C previously held a formula (=B*1000, shared across C3:C23).
# Assigning a literal .Value below replaces the formula with a plain number,
# matching the diff (formulas removed, static values written instead).

# Row 2
$ws.Range("B2").Value = 0.1
$ws.Range("C2").Value = 3

# Row 3
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 25

# Row 4
$ws.Range("B4").Value = 0.0001
$ws.Range("C4").Value = 0.1

# Row 5
$ws.Range("B5").Value = 0.0001
$ws.Range("C5").Value = 0.1

# Row 6
$ws.Range("B6").Value = 0.01
$ws.Range("C6").Value = 0.05

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 100

# Row 8 (B8 unchanged)
$ws.Range("C8").Value = 5

# Row 9 (B9 unchanged; C9 formula -> static, same value)
$ws.Range("C9").Value = 0

# Row 10 (B10 unchanged; C10 formula -> static, same value)
$ws.Range("C10").Value = 0

# Row 11
$ws.Range("B11").Value = 0.001
$ws.Range("C11").Value = 100

# Row 12 (B12 unchanged)
$ws.Range("C12").Value = 0.01

# Row 13 (B13 unchanged; C13 formula -> static, same value)
$ws.Range("C13").Value = 0

# Row 14 (B14 unchanged; C14 formula -> static, same value)
$ws.Range("C14").Value = 0

# Row 15
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 10

# Row 16 (B16 unchanged; C16 formula -> static, same value)
$ws.Range("C16").Value = 0.01

# Row 17
$ws.Range("B17").Value = 0.00001
$ws.Range("C17").Value = 0.001

# Row 18: no change (already static values)

# Row 19
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 10

# Row 20
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 25

# Row 21
$ws.Range("B21").Value = 0.001
$ws.Range("C21").Value = 0.01

# Row 22
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 10

# Row 23 (B23 unchanged)
$ws.Range("C23").Value = 3

# --- Update sheet selection -------------------------------------------------
$ws.Range("C8").Select()

# --- Update workbook window geometry (position/size of the saved window) ---
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 760
$win.Width = 19820
$win.Height = 16840
